$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 15.70818033333333
$ws.Range("N2").Value = 47.12454099999999
$ws.Range("O2").Value = 0.3220467100482788
$ws.Range("P2").Value = 0.334408980496766
$ws.Range("Q2").Value = 4.711789120365888
$ws.Range("R2").Value = 42.40610208329299
$ws.Range("S2").Value = 0.3220467100482788
$ws.Range("T2").Value = 0.334408980496766

# Row 3
$ws.Range("O3").Value = 0.1200026410479322
$ws.Range("P3").Value = 0.1246091315254933
$ws.Range("S3").Value = 0.1200026410479322
$ws.Range("T3").Value = 0.1246091315254933

# Row 4
$ws.Range("M4").Value = 10.959131
$ws.Range("N4").Value = 32.877393
$ws.Range("O4").Value = 0.2246824271585863
$ws.Range("P4").Value = 0.2333072161810874
$ws.Range("Q4").Value = 3.287275363454334
$ws.Range("R4").Value = 29.58547827108901
$ws.Range("S4").Value = 0.2246824271585863
$ws.Range("T4").Value = 0.2333072161810874

# Row 5
$ws.Range("M5").Value = 5.4093935
$ws.Range("N5").Value = 10.818787
$ws.Range("O5").Value = 0.1109025579706895
$ws.Range("P5").Value = 0.07677315161290731
$ws.Range("Q5").Value = 1.622589052341833
$ws.Range("R5").Value = 9.735534314051002
$ws.Range("S5").Value = 0.1109025579706895
$ws.Range("T5").Value = 0.07677315161290731

# Row 6
$ws.Range("M6").Value = 10.84612833333333
$ws.Range("N6").Value = 32.538385
$ws.Range("O6").Value = 0.2223656637745133
$ws.Range("P6").Value = 0.230901520183746
$ws.Range("Q6").Value = 3.253379347233889
$ws.Range("R6").Value = 29.280414125105
$ws.Range("S6").Value = 0.2223656637745133
$ws.Range("T6").Value = 0.230901520183746
